# Update "想去人数" (want-to-go count, column F) figures across the four
# sheets (展览 / 演出 / 本地生活 / 全部类型) to match the refreshed
# gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 3247
$ws.Range("F7").Value = 4732
$ws.Range("F8").Value = 447
$ws.Range("F9").Value = 271
$ws.Range("F13").Value = 16
$ws.Range("F15").Value = 643
$ws.Range("F20").Value = 332
$ws.Range("F21").Value = 4693
$ws.Range("F25").Value = 5840
$ws.Range("F27").Value = 1180
$ws.Range("F28").Value = 236
$ws.Range("F29").Value = 652
$ws.Range("F31").Value = 3
$ws.Range("F32").Value = 75
$ws.Range("F33").Value = 117
$ws.Range("F34").Value = 813
$ws.Range("F36").Value = 737
$ws.Range("F37").Value = 756

# 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 10

# 本地生活 (Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 1077

# 全部类型 (All types) — union of the other three sheets
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 1077
$ws.Range("F9").Value = 3247
$ws.Range("F11").Value = 4732
$ws.Range("F12").Value = 447
$ws.Range("F13").Value = 271
$ws.Range("F17").Value = 16
$ws.Range("F19").Value = 643
$ws.Range("F25").Value = 332
$ws.Range("F26").Value = 4693
$ws.Range("F30").Value = 5840
$ws.Range("F32").Value = 1180
$ws.Range("F33").Value = 236
$ws.Range("F34").Value = 652
$ws.Range("F36").Value = 3
$ws.Range("F37").Value = 10
$ws.Range("F38").Value = 75
$ws.Range("F39").Value = 117
$ws.Range("F40").Value = 813
$ws.Range("F42").Value = 737
$ws.Range("F43").Value = 756
